$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The leaderboard now shows the last question answered "ok" for the two
# rows that previously had a blank Status column (row 6 "Paso, paso..."
# and row 9 "Mostrar ultima pregunta contestada ok").
$ws.Range("C6").Value = "ok"
$ws.Range("C9").Value = "ok"

# The author's selection moved from C6 to C10 after making the edit.
$ws.Range("C10").Select()
